# Update gh-pages to output generated at 456a3b4
# Apply updated "想去人数" (F) and "最低票价" (G) values for the 展览 (Exhibition)
# sheet, and mirror the same underlying row updates into the 全部类型 (All types)
# sheet, which aggregates rows from all the other sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F2").Value = 14162
$ws1.Range("G2").Value = 95

$ws1.Range("F6").Value = 1054

$ws1.Range("F7").Value = 13963

$ws1.Range("F8").Value = 15105

$ws1.Range("F10").Value = 25

$ws1.Range("F21").Value = 1172

$ws1.Range("F24").Value = 5865

$ws1.Range("F26").Value = 1072

$ws1.Range("F27").Value = 5477

$ws1.Range("F31").Value = 365

# --- Sheet "全部类型" (same records, offset by the interleaved rows from the
#     other category sheets) ---
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F2").Value = 14162
$ws4.Range("G2").Value = 95

$ws4.Range("F7").Value = 1054

$ws4.Range("F8").Value = 13963

$ws4.Range("F9").Value = 15105

$ws4.Range("F11").Value = 25

$ws4.Range("F22").Value = 1172

$ws4.Range("F26").Value = 5865

$ws4.Range("F28").Value = 1072

$ws4.Range("F29").Value = 5477

$ws4.Range("F33").Value = 365
